# Adv Filter changes for Item, ItemClass, ItemType and Identifier Type screen
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Login sheet: update stored username value
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item("Login")
$login.Range("A2").Value = "balasanthanam92@gmail.com"

# ---------------------------------------------------------------------------
# 2) Create the four new Advanced-Filter worksheets, appended at the end of
#    the workbook (after the existing "FacilityAdvancedFilter" sheet).
# ---------------------------------------------------------------------------

# --- ItemclassAdvFilter ----------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsItemClassFilter = $wb.Worksheets.Add($null, $last)
$wsItemClassFilter.Name = "ItemclassAdvFilter"
$wsItemClassFilter.Range("A1").Value = "sCode"
$wsItemClassFilter.Range("B1").Value = "sName"
$wsItemClassFilter.Range("C1").Value = "sDescription"
$wsItemClassFilter.Range("A2").Value = "test8"
$wsItemClassFilter.Range("B2").Value = "test8"
$wsItemClassFilter.Range("C2").Value = "test8"

# --- ItemTypeAdvFilter ------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsItemTypeFilter = $wb.Worksheets.Add($null, $last)
$wsItemTypeFilter.Name = "ItemTypeAdvFilter"
$wsItemTypeFilter.Range("A1").Value = "sCode"
$wsItemTypeFilter.Range("B1").Value = "sName"
$wsItemTypeFilter.Range("C1").Value = "sDescription"
$wsItemTypeFilter.Range("A2").Value = "CLASS1"
$wsItemTypeFilter.Range("B2").Value = "CLASS1"
$wsItemTypeFilter.Range("C2").Value = "CLASS1"

# --- IdentifierTypeAdvFilter -------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsIdentifierTypeFilter = $wb.Worksheets.Add($null, $last)
$wsIdentifierTypeFilter.Name = "IdentifierTypeAdvFilter"
$wsIdentifierTypeFilter.Range("A1").Value = "sCode"
$wsIdentifierTypeFilter.Range("B1").Value = "sName"
$wsIdentifierTypeFilter.Range("C1").Value = "sIsRFID"
$wsIdentifierTypeFilter.Range("D1").Value = "sProtocol"
$wsIdentifierTypeFilter.Range("E1").Value = "sFromDate"
$wsIdentifierTypeFilter.Range("F1").Value = "sToDate"
$wsIdentifierTypeFilter.Range("A2").Value = "RFID"
$wsIdentifierTypeFilter.Range("B2").Value = "RFID"
$wsIdentifierTypeFilter.Range("C2").Value = "Yes"
$wsIdentifierTypeFilter.Range("D2").Value = "RFCODE"
# Dates are stored as literal text in the source workbook, not as Excel
# date serials - force the text number format before writing the values.
$wsIdentifierTypeFilter.Range("E2").NumberFormat = "@"
$wsIdentifierTypeFilter.Range("E2").Value = "12-12-2016"
$wsIdentifierTypeFilter.Range("F2").NumberFormat = "@"
$wsIdentifierTypeFilter.Range("F2").Value = "15-11-2017"

# --- ItemAdvancedFilter ------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsItemFilter = $wb.Worksheets.Add($null, $last)
$wsItemFilter.Name = "ItemAdvancedFilter"
$wsItemFilter.Range("A1").Value = "sItem"
$wsItemFilter.Range("B1").Value = "sName"
$wsItemFilter.Range("C1").Value = "sItemClass"
$wsItemFilter.Range("D1").Value = "sItemType"
$wsItemFilter.Range("E1").Value = "sManufacture"
$wsItemFilter.Range("F1").Value = "sManufactureItem"
$wsItemFilter.Range("G1").Value = "sInventory"
$wsItemFilter.Range("H1").Value = "sInventoryType"

$wsItemFilter.Range("A2").Value = "CSVITEMCODE100"
$wsItemFilter.Range("B2").Value = "CSVITEMNAME100"
$wsItemFilter.Range("B2").WrapText = $true
$wsItemFilter.Range("C2").Value = "VEHICLE"
$wsItemFilter.Range("C2").WrapText = $true
$wsItemFilter.Range("D2").Value = "CLASS1"
$wsItemFilter.Range("E2").Value = "Toyota"
$wsItemFilter.Range("F2").Value = "T100"
$wsItemFilter.Range("G2").Value = "No"
$wsItemFilter.Range("H2").Value = "Non-Serialized"

# ---------------------------------------------------------------------------
# 3) Selection / active-sheet bookkeeping to mirror the recorded state:
#    - "Itemclass" remembers A1 as its last selection
#    - "FacilityAdvancedFilter" stays on K10 but is no longer the active tab
#    - "Login" becomes the active sheet/tab, with D13 selected
# ---------------------------------------------------------------------------
$itemclass = $wb.Worksheets.Item("Itemclass")
[void]$itemclass.Activate()
[void]$itemclass.Range("A1").Select()

$facilityFilter = $wb.Worksheets.Item("FacilityAdvancedFilter")
[void]$facilityFilter.Activate()
[void]$facilityFilter.Range("K10").Select()
# Normalise the stray duplicate-font style that used to be on I2 so it
# matches the plain default formatting used elsewhere in the row.
$facilityFilter.Range("I2").Font.Name = "Calibri"
$facilityFilter.Range("I2").Font.Size = 11
$facilityFilter.Range("I2").Font.Color = 0

[void]$login.Activate()
[void]$login.Range("D13").Select()
